$wb = $excel.ActiveWorkbook

# --- Rename "Requirements" to "SoftwareRequirements" ---
$reqSheet = $wb.Worksheets.Item("Requirements")
$reqSheet.Name = "SoftwareRequirements"

# --- Insert a new blank sheet "Interface Def and Reqs" right after it ---
$newSheet = $wb.Worksheets.Add($null, $reqSheet)
$newSheet.Name = "Interface Def and Reqs"

# --- Fix up the A9:A10 cell styling on SoftwareRequirements sheet: it had a
#     stray font-applied style (index 11) that should instead match the
#     plain bordered style used by the rest of the column (same as A8). ---
$srcStyle = $reqSheet.Range("A8")
$dstStyle = $reqSheet.Range("A9:A10")
$srcStyle.Copy()
$dstStyle.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Activate "Requiremnts Unknowns" (now the 3rd tab) so it becomes the
#     selected/active sheet when the workbook is reopened. Re-fetch the
#     reference since the collection changed after Add(). ---
$unkSheet = $wb.Worksheets.Item("Requiremnts Unknowns")
$unkSheet.Activate()
